# Generate Report for Handoff
# The "d3775bf9-f54e-4da7-87aa-79d1c9bb0ad8" file is now ready for handoff again:
#  - Overview sheet: Status columns (zh-cn / de-de) -> "Ready for handoff",
#    and the Latest Handoff Date -> the new handoff timestamp.
#  - zh-cn sheet: Status -> "Ready for handoff", Latest Handoff Datetime -> new timestamp.
#  - de-de sheet: Status -> "Ready for handoff", Latest Handoff Datetime -> new timestamp.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status
$wsOverview.Range("D3").Value = "2016-37-19 20:37:53"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("E3").Value = "2016-03-19 20:37:51"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("E3").Value = "2016-03-19 20:37:53"
